{"js": "// Update the date stamp and all of the two-digit-by-two-digit\n// multiplication problems in the practice sheet to the new values.\n//\n// Each \"old\" string below occurs exactly once in the document body, so a\n// simple text search + replace per pair is unambiguous and order-safe\n// (no \"new\" value ever collides with a not-yet-processed \"old\" value).\nconst replacements = [\n  [\"2025-02-26 Wednesday\", \"2025-02-27 Thursday\"],\n  [\"26\u00d784=\", \"83\u00d736=\"],\n  [\"21\u00d726=\", \"39\u00d795=\"],\n  [\"91\u00d717=\", \"60\u00d751=\"],\n  [\"19\u00d774=\", \"83\u00d723=\"],\n  [\"42\u00d712=\", \"94\u00d731=\"],\n  [\"47\u00d795=\", \"29\u00d740=\"],\n  [\"11\u00d792=\", \"45\u00d744=\"],\n  [\"80\u00d762=\", \"34\u00d795=\"],\n  [\"92\u00d732=\", \"22\u00d770=\"],\n  [\"30\u00d720=\", \"14\u00d754=\"],\n  [\"82\u00d751=\", \"77\u00d727=\"],\n  [\"20\u00d734=\", \"25\u00d772=\"],\n  [\"79\u00d780=\", \"26\u00d750=\"],\n  [\"12\u00d741=\", \"88\u00d745=\"],\n  [\"91\u00d761=\", \"54\u00d754=\"],\n  [\"65\u00d785=\", \"20\u00d750=\"],\n  [\"43\u00d729=\", \"13\u00d787=\"],\n  [\"30\u00d786=\", \"34\u00d795=\"],\n  [\"61\u00d758=\", \"76\u00d722=\"],\n  [\"75\u00d758=\", \"63\u00d742=\"],\n  [\"83\u00d773=\", \"46\u00d749=\"],\n  [\"94\u00d720=\", \"78\u00d719=\"],\n  [\"39\u00d769=\", \"58\u00d732=\"],\n  [\"70\u00d752=\", \"63\u00d795=\"],\n  [\"14\u00d789=\", \"56\u00d741=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date stamp and all of the two-digit-by-two-digit\n# multiplication problems in the practice sheet to the new values.\n#\n# Each \"old\" string below occurs exactly once in the document, so a\n# simple Find/Replace per pair (restricted to exact whole-string matches)\n# is unambiguous and order-safe (no \"new\" value ever collides with a\n# not-yet-processed \"old\" value).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2025-02-26 Wednesday\", \"2025-02-27 Thursday\")\n    ,@(\"26\u00d784=\", \"83\u00d736=\")\n    ,@(\"21\u00d726=\", \"39\u00d795=\")\n    ,@(\"91\u00d717=\", \"60\u00d751=\")\n    ,@(\"19\u00d774=\", \"83\u00d723=\")\n    ,@(\"42\u00d712=\", \"94\u00d731=\")\n    ,@(\"47\u00d795=\", \"29\u00d740=\")\n    ,@(\"11\u00d792=\", \"45\u00d744=\")\n    ,@(\"80\u00d762=\", \"34\u00d795=\")\n    ,@(\"92\u00d732=\", \"22\u00d770=\")\n    ,@(\"30\u00d720=\", \"14\u00d754=\")\n    ,@(\"82\u00d751=\", \"77\u00d727=\")\n    ,@(\"20\u00d734=\", \"25\u00d772=\")\n    ,@(\"79\u00d780=\", \"26\u00d750=\")\n    ,@(\"12\u00d741=\", \"88\u00d745=\")\n    ,@(\"91\u00d761=\", \"54\u00d754=\")\n    ,@(\"65\u00d785=\", \"20\u00d750=\")\n    ,@(\"43\u00d729=\", \"13\u00d787=\")\n    ,@(\"30\u00d786=\", \"34\u00d795=\")\n    ,@(\"61\u00d758=\", \"76\u00d722=\")\n    ,@(\"75\u00d758=\", \"63\u00d742=\")\n    ,@(\"83\u00d773=\", \"46\u00d749=\")\n    ,@(\"94\u00d720=\", \"78\u00d719=\")\n    ,@(\"39\u00d769=\", \"58\u00d732=\")\n    ,@(\"70\u00d752=\", \"63\u00d795=\")\n    ,@(\"14\u00d789=\", \"56\u00d741=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
